$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.077005770509684
$ws.Range("D2").Value = 1.07910128515306
$ws.Range("E2").Value = 1.080506282421036
$ws.Range("F2").Value = 1.090788279320382
$ws.Range("I2").Value = 1.057794804003207
$ws.Range("J2").Value = 1.081902390064228
$ws.Range("K2").Value = 1.081778218676892
$ws.Range("L2").Value = 1.083179538724866
$ws.Range("M2").Value = 1.093434936186371
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.078262152919711
$ws.Range("D3").Value = 1.080123569077601
$ws.Range("E3").Value = 1.0816196056051
$ws.Range("F3").Value = 1.091963856306752
$ws.Range("I3").Value = 1.058189823679921
$ws.Range("J3").Value = 1.082817381373135
$ws.Range("K3").Value = 1.08261801192799
$ws.Range("L3").Value = 1.084110416354603
$ws.Range("M3").Value = 1.094429845571146
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.079074988713257
$ws.Range("D4").Value = 1.080784886679979
$ws.Range("E4").Value = 1.08234011449908
$ws.Range("F4").Value = 1.092724739205675
$ws.Range("I4").Value = 1.058444131242422
$ws.Range("J4").Value = 1.083408752860226
$ws.Range("K4").Value = 1.083160632370092
$ws.Range("L4").Value = 1.084712274667444
$ws.Range("M4").Value = 1.095073238099369
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.07941667629054
$ws.Range("D5").Value = 1.081062865215782
$ws.Range("E5").Value = 1.082643044904911
$ws.Range("F5").Value = 1.093044664867024
$ws.Range("I5").Value = 1.058550732271079
$ws.Range("J5").Value = 1.083657201620814
$ws.Range("K5").Value = 1.083388563778725
$ws.Range("L5").Value = 1.084965181618333
$ws.Range("M5").Value = 1.095343630351736
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.079474045518925
$ws.Range("D6").Value = 1.081109536753047
$ws.Range("E6").Value = 1.08269390995406
$ws.Range("F6").Value = 1.09309838478838
$ws.Range("I6").Value = 1.058568612897437
$ws.Range("J6").Value = 1.083698907688429
$ws.Range("K6").Value = 1.083426823573181
$ws.Range("L6").Value = 1.08500763913712
$ws.Range("M6").Value = 1.095389025156864
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.079079554469259
$ws.Range("D7").Value = 1.080788601196348
$ws.Range("E7").Value = 1.082344162156765
$ws.Range("F7").Value = 1.092729013868146
$ws.Range("I7").Value = 1.058445556867493
$ws.Range("J7").Value = 1.083412073286951
$ws.Range("K7").Value = 1.083163678731247
$ws.Range("L7").Value = 1.084715654470241
$ws.Range("M7").Value = 1.095076851446013
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.077430398438902
$ws.Range("D8").Value = 1.079446806240233
$ws.Range("E8").Value = 1.08088251218453
$ws.Range("F8").Value = 1.091185529127079
$ws.Range("I8").Value = 1.057928571539647
$ws.Range("J8").Value = 1.082211759032395
$ws.Range("K8").Value = 1.082062193231961
$ws.Range("L8").Value = 1.083494233562162
$ws.Range("M8").Value = 1.093771249888957
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.074523309072561
$ws.Range("D9").Value = 1.077081059230145
$ws.Range("E9").Value = 1.078307724692235
$ws.Range("F9").Value = 1.088467236221835
$ws.Range("I9").Value = 1.057007622998563
$ws.Range("J9").Value = 1.080091322565266
$ws.Range("K9").Value = 1.080115202330392
$ws.Range("L9").Value = 1.081338196811558
$ws.Range("M9").Value = 1.0914676505683
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.072584407486494
$ws.Range("D10").Value = 1.075502936145591
$ws.Range("E10").Value = 1.07659168394266
$ws.Range("F10").Value = 1.086655988928857
$ws.Range("I10").Value = 1.056386929655789
$ws.Range("J10").Value = 1.078674040325063
$ws.Range("K10").Value = 1.078813088750745
$ws.Range("L10").Value = 1.079898257941541
$ws.Range("M10").Value = 1.089929854420749
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.071744613369308
$ws.Range("D11").Value = 1.074819349437818
$ws.Range("E11").Value = 1.075848717499336
$ws.Range("F11").Value = 1.085871905747681
$ws.Range("I11").Value = 1.056116558694187
$ws.Range("J11").Value = 1.078059456532414
$ws.Range("K11").Value = 1.078248267235919
$ws.Range("L11").Value = 1.079274120380863
$ws.Range("M11").Value = 1.089263465691169
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.071432638285628
$ws.Range("D12").Value = 1.074565396125411
$ws.Range("E12").Value = 1.075572758744415
$ws.Range("F12").Value = 1.08558069057062
$ws.Range("I12").Value = 1.056015888729518
$ws.Range("J12").Value = 1.077831037263787
$ws.Range("K12").Value = 1.078038316106229
$ws.Range("L12").Value = 1.07904219110331
$ws.Range("M12").Value = 1.089015860665531
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.071499559767912
$ws.Range("D13").Value = 1.074619871749454
$ws.Range("E13").Value = 1.07563195231298
$ws.Range("F13").Value = 1.085643155963778
$ws.Range("I13").Value = 1.056037493748902
$ws.Range("J13").Value = 1.077880040128721
$ws.Range("K13").Value = 1.078083358203589
$ws.Range("L13").Value = 1.079091945124542
$ws.Range("M13").Value = 1.089068976364591
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.07171882620477
$ws.Range("D14").Value = 1.074798358376637
$ws.Range("E14").Value = 1.0758259064254
$ws.Range("F14").Value = 1.085847833244002
$ws.Range("I14").Value = 1.056108242227772
$ws.Range("J14").Value = 1.078040578097425
$ws.Range("K14").Value = 1.078230915703852
$ws.Range("L14").Value = 1.079254951017571
$ws.Range("M14").Value = 1.08924300019439
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.071853918360469
$ws.Range("D15").Value = 1.074908324699879
$ws.Range("E15").Value = 1.075945409507723
$ws.Range("F15").Value = 1.085973945362668
$ws.Range("I15").Value = 1.056151800590155
$ws.Range("J15").Value = 1.078139472837962
$ws.Range("K15").Value = 1.07832181066393
$ws.Range("L15").Value = 1.079355371457559
$ws.Range("M15").Value = 1.089350211560925
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.072640136645868
$ws.Range("D16").Value = 1.075548298181667
$ws.Range("E16").Value = 1.076640993954704
$ws.Range("F16").Value = 1.086708029996296
$ws.Range("I16").Value = 1.05640483937007
$ws.Range("J16").Value = 1.078714809314755
$ws.Range("K16").Value = 1.078850552956252
$ws.Range("L16").Value = 1.07993966640196
$ws.Range("M16").Value = 1.089974069510301
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.073133245363373
$ws.Range("D17").Value = 1.075949669013213
$ws.Range("E17").Value = 1.077077338852386
$ws.Range("F17").Value = 1.087168553967981
$ws.Range("I17").Value = 1.056563133221976
$ws.Range("J17").Value = 1.079075462959168
$ws.Range("K17").Value = 1.07918195093558
$ws.Range("L17").Value = 1.080306008267181
$ws.Range("M17").Value = 1.090365260735188
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.073420844648949
$ws.Range("D18").Value = 1.076183757897347
$ws.Range("E18").Value = 1.077331860288496
$ws.Range("F18").Value = 1.087437189107983
$ws.Range("I18").Value = 1.056655308382049
$ws.Range("J18").Value = 1.079285740245738
$ws.Range("K18").Value = 1.079375153596188
$ws.Range("L18").Value = 1.080519628156727
$ws.Range("M18").Value = 1.090593386580334
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.073518904801964
$ws.Range("D19").Value = 1.076263572139001
$ws.Range("E19").Value = 1.077418647066124
$ws.Range("F19").Value = 1.087528790071242
$ws.Range("I19").Value = 1.056686711467831
$ws.Range("J19").Value = 1.079357424846791
$ws.Range("K19").Value = 1.07944101443892
$ws.Range("L19").Value = 1.080592456695122
$ws.Range("M19").Value = 1.090671163268689
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.073080341831044
$ws.Range("D20").Value = 1.075906608207653
$ws.Range("E20").Value = 1.077030522272196
$ws.Range("F20").Value = 1.087119142120986
$ws.Range("I20").Value = 1.056546165824635
$ws.Range("J20").Value = 1.07903677713194
$ws.Range("K20").Value = 1.07914640502721
$ws.Range("L20").Value = 1.080266709585771
$ws.Range("M20").Value = 1.090323294738373
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.071654258778074
$ws.Range("D21").Value = 1.074745799594871
$ws.Range("E21").Value = 1.07576879145285
$ws.Range("F21").Value = 1.085787560120446
$ws.Range("I21").Value = 1.056087415254854
$ws.Range("J21").Value = 1.077993307419713
$ws.Range("K21").Value = 1.078187467884625
$ws.Range("L21").Value = 1.079206952528779
$ws.Range("M21").Value = 1.089191756695459
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.070757401728996
$ws.Range("D22").Value = 1.074015728144503
$ws.Range("E22").Value = 1.074975560235841
$ws.Range("F22").Value = 1.084950504165096
$ws.Range("I22").Value = 1.055797579169668
$ws.Range("J22").Value = 1.077336451942168
$ws.Range("K22").Value = 1.077583670158446
$ws.Range("L22").Value = 1.078540080434742
$ws.Range("M22").Value = 1.088479858414581
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.071232864521775
$ws.Range("D23").Value = 1.074402774617867
$ws.Range("E23").Value = 1.075396061084092
$ws.Range("F23").Value = 1.085394228512842
$ws.Range("I23").Value = 1.05595135985575
$ws.Range("J23").Value = 1.077684738365737
$ws.Range("K23").Value = 1.077903838239686
$ws.Range("L23").Value = 1.078893655594869
$ws.Range("M23").Value = 1.088857292841862
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.073104246730219
$ws.Range("D24").Value = 1.075926065606205
$ws.Range("E24").Value = 1.077051676642378
$ws.Range("F24").Value = 1.087141469148211
$ws.Range("I24").Value = 1.056553833139997
$ws.Range("J24").Value = 1.079054257857476
$ws.Range("K24").Value = 1.079162466990378
$ws.Range("L24").Value = 1.080284467157242
$ws.Range("M24").Value = 1.09034225751858
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.075275000564557
$ws.Range("D25").Value = 1.0776928258892
$ws.Range("E25").Value = 1.078973278376621
$ws.Range("F25").Value = 1.089169807356515
$ws.Range("I25").Value = 1.05724689350608
$ws.Range("J25").Value = 1.080640145467534
$ws.Range("K25").Value = 1.080619267594971
$ws.Range("L25").Value = 1.081896033891163
$ws.Range("M25").Value = 1.092063544631949
